$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "Jett:"
$ws.Range("B17").Value = '${jett:cellRef(16, 1)}'
$ws.Range("D17").Value = '${jett:cellRef(16, 1, 2, 3)}'

$ws.Range("A18").Value = "Static:"
$ws.Range("B18").Value = '${java.lang.String.format(''%s supports static method calling!'', testBean2)}'
